# "Add files via upload" - appends five newly-scraped credential rows
# (Username / Password) to the bottom of the Outlook-credentials table on
# Sheet1. Column A gets a mailto: hyperlink whose display text is the
# e-mail address (same convention as every existing row); column B gets
# the plain-text password.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("gvarderesyan@mail.ru", "garnik1986"),
    @("hov.man1998@mail.ru", "heghnar1969"),
    @("dianamartayan@bk.ru", "hripsime1952"),
    @("lauraasatryan@gmail.com", "armen1966"),
    @("arpine_mxitaryan@mail.ru", "arpine1979")
)

# Existing row used as a formatting template, so the new rows pick up the
# same cell styles as the rest of the table instead of minting new ones.
$templateA = $ws.Cells.Item(3, 1)

$startRow = 17
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $email = $newData[$i][0]
    $password = $newData[$i][1]

    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)

    $cellA.Value = $email
    $ws.Hyperlinks.Add($cellA, "mailto:$email") | Out-Null

    # Hyperlinks.Add stamps its own built-in "Hyperlink" style; restore the
    # worksheet's own hyperlink look (left-aligned, themed) by copying the
    # format from an existing username cell.
    $templateA.Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null

    $cellB.Value = $password
}

$excel.CutCopyMode = $false

# Leave the selection where the author's cursor ended up after entering the
# last row of data.
$ws.Range("B22").Select() | Out-Null
